{"js": "// Bold + color (\"2C3E50\") the quantitative metrics (percentages, dollar\n// amounts, and large numbers) inside specific bullet/body paragraphs, as\n// described by the diff. We scope each substring search to the exact\n// paragraph it belongs to (by index) so we never touch look-alike numbers\n// that appear elsewhere in the document (e.g. the \"23% to 64%\" phrase that\n// also shows up in the PROFESSIONAL SUMMARY and KEY PROJECTS sections,\n// which the diff leaves untouched).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Map of paragraph index -> list of exact substrings (in left-to-right\n// order) that must become bold + colored (#2C3E50) in that paragraph.\nconst plan = [\n  { index: 8, terms: [\"23%\", \"64%\"] },\n  { index: 10, terms: [\"87%\", \"71%\", \"\\u00B14.2%\", \"\\u00B12.1%\"] },\n  { index: 30, terms: [\"1,200\"] },\n  { index: 45, terms: [\"$400M\", \"$1B\"] },\n  { index: 62, terms: [\"73.5%\", \"$4.7M\"] },\n  { index: 64, terms: [\"87%\", \"71%\"] },\n];\n\nfor (const { index, terms } of plan) {\n  const para = items[index];\n  for (const term of terms) {\n    const found = para.search(term, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n\n    for (const range of found.items) {\n      range.font.bold = true;\n      range.font.color = \"2C3E50\";\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Bold + color (#2C3E50) the quantitative metrics (percentages, dollar\n# amounts, and large numbers) inside specific bullet/body paragraphs, as\n# described by the diff. Each substring search is scoped to the exact\n# paragraph it belongs to (by 1-based COM Paragraphs index) so we never\n# touch look-alike numbers that appear elsewhere in the document (e.g. the\n# \"23% to 64%\" phrase that also shows up in the PROFESSIONAL SUMMARY and\n# KEY PROJECTS sections, which the diff leaves untouched).\n\n$doc = $word.ActiveDocument\n\n# #2C3E50 as an OLE/VBA BGR-packed color long: R=0x2C G=0x3E B=0x50\n# -> B*65536 + G*256 + R\n$metricColor = 0x50 * 65536 + 0x3E * 256 + 0x2C\n\n# Map of 1-based Paragraphs() index -> list of exact substrings (in\n# left-to-right order) that must become bold + colored in that paragraph.\n$plan = @(\n    @{ Index = 9;  Terms = @(\"23%\", \"64%\") },\n    @{ Index = 11; Terms = @(\"87%\", \"71%\", \"\u00b14.2%\", \"\u00b12.1%\") },\n    @{ Index = 31; Terms = @(\"1,200\") },\n    @{ Index = 46; Terms = @(\"$400M\", \"$1B\") },\n    @{ Index = 63; Terms = @(\"73.5%\", \"$4.7M\") },\n    @{ Index = 65; Terms = @(\"87%\", \"71%\") }\n)\n\nforeach ($entry in $plan) {\n    $paraRange = $doc.Paragraphs($entry.Index).Range\n    $paraStart = $paraRange.Start\n    $paraEnd = $paraRange.End\n\n    foreach ($term in $entry.Terms) {\n        $searchRange = $doc.Range($paraStart, $paraEnd)\n        $found = $searchRange.Find.Execute($term, $true)\n        if ($found) {\n            $searchRange.Font.Bold = 1\n            $searchRange.Font.Color = $metricColor\n        }\n    }\n}\n"}
